$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "CubeA"

# Add a new row 16 with averaged intensity values computed via the
# Gaussian Quadrature scheme (HKL index 14, same label as row 15)
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "HexGrid-60degTilt5degRes"
$ws.Cells.Item(16, 3).Value = 0.9971603071729759
$ws.Cells.Item(16, 4).Value = 1.012708961054975
$ws.Cells.Item(16, 5).Value = 0.9952648695185737
$ws.Cells.Item(16, 6).Value = 1.001633974110304
$ws.Cells.Item(16, 7).Value = 0.9971603071729759
$ws.Cells.Item(16, 8).Value = 1.012708961054975
$ws.Cells.Item(16, 9).Value = 0.9979085223676277
$ws.Cells.Item(16, 10).Value = 1.007545575179965
$ws.Cells.Item(16, 11).Value = 0.998206045989162
$ws.Cells.Item(16, 12).Value = 1.011503242338255
$ws.Cells.Item(16, 13).Value = 0.9971603071729759
$ws.Cells.Item(16, 14).Value = 1.003986915286774
$ws.Cells.Item(16, 15).Value = 1.001692027964207
$ws.Cells.Item(16, 16).Value = 1.00274143721648

# Match the formatting of the preceding row's HKL index cell (A15),
# which is bold, bordered, and centered
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
